# The "Dataset Details" slide (slide 5) has its body placeholder
# (the second shape in the shape tree) nudged/resized: the left edge
# moves in by 1 EMU and the box is widened from 5640900 EMU to
# 6298349 EMU, while the top/height stay the same.
#
# PowerPoint's COM object model works in points (1 pt = 12700 EMU), so
# the literals below are the point values that round-trip to the exact
# target EMU figures (399816 / 6298349) used by the canonical OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)

$shp.Left = 31.4816
$shp.Top = 101.99252
$shp.Width = 495.933
$shp.Height = 260.7178
